$wb = $excel.ActiveWorkbook

# --- "unknowns" sheet: just a view-state change (selection moved) ---
$unknowns = $wb.Worksheets.Item("unknowns")
$unknowns.Activate()
$unknowns.Range("E21").Select()

# --- "data" sheet: just a view-state change (selection moved) ---
$data = $wb.Worksheets.Item("data")
$data.Activate()
$data.Range("E40").Select()

# --- "system" sheet: add new outline steps 8-10 about the mass spec model,
#     Faraday/IonCounter detection, and future interference parameters ---
$system = $wb.Worksheets.Item("system")
$system.Activate()

$system.Range("A39").Value = "8. From the mass spectrometer model, determine whether each collector is a Faraday or Ion Counter"

$system.Range("B41").Value = 'isa(massSpec.collectorArray(["Ax"],:).collectors, "faradayModel")'
$system.Range("H41").Value = $true
$system.Range("B42").Value = 'isa(massSpec.collectorArray(["H1"],:).collectors, "faradayModel")'
$system.Range("H42").Value = $true

$system.Range("A44").Value = "9. If any Faradays: "
$system.Range("A47").Value = "9. If any IonCounters: "
$system.Range("B45").Value = "enable betaFaraday, upMassTailFaraday, and downMassTailFaraday as model parameters"
$system.Range("B48").Value = "enable betaIonCounter, darkNoise, upMassIonCounter, and downMassIonCounter as model parameters"

$system.Range("A50").Value = "10. Future: add interferences"
$system.Range("B52").Value = "parameterization TBA"

$system.Range("B53").Select()

# --- new "syndata" sheet, appended after "system", outlining synthetic-data inputs ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$syndata = $wb.Worksheets.Add($null, $lastSheet)
$syndata.Name = "syndata"

$syndata.Range("A1").Value = "To create synthetic data, we need to combine data from the method, the mass spectrometer, and from user-specific input"

$syndata.Range("A3").Value = "To create the true modelParameters:"
$syndata.Range("A4").Value = 'Initiate the system as described in the "system" sheet.'

$syndata.Range("A6").Value = "uses: method, massSpec, user input"

$syndata.Range("A8").Value = 'Additional information needed to produce the data from the "data" sheet:'

$syndata.Range("B10").Value = "Intensity function (as spline knots values)"
$syndata.Range("B11").Value = "Measurement time stamps (need integration times from method for consistency)"

$syndata.Range("A10").Select()
